# Commit: "update test case ids for store pets"
#
# - Sheet1 ("id"/"petId" columns A & B, rows 2-11): renumber the test ids
#   from the old 124-133 / 167-176 ranges down to a clean 1-10 sequence.
# - Sheet2 (A2): update the single order id from 7 to 2.
# - The active sheet/tab flips from Sheet2 back to Sheet1, and each sheet's
#   remembered selection moves accordingly (Sheet1 -> H7, Sheet2 -> D9).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: renumber id (col A) and petId (col B) for rows 2..11 to 1..10 ---
for ($i = 0; $i -lt 10; $i++) {
    $row = 2 + $i
    $newValue = $i + 1
    $ws1.Cells.Item($row, 1).Value = $newValue
    $ws1.Cells.Item($row, 2).Value = $newValue
}

# --- Sheet2: orderId in A2 goes from 7 to 2 ---
$ws2.Range("A2").Value = 2

# --- Make Sheet1 the active sheet/tab again, with its remembered selection ---
$ws1.Activate()
$ws1.Range("H7").Select()

# --- Sheet2 keeps its own remembered selection for when it's revisited ---
$ws2.Range("D9").Select()

# Re-activate Sheet1 last so it is the sheet that is actually active/visible
$ws1.Activate()
